$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("751:751").Insert()

$ws.Range("A751").Value = 3
$ws.Range("B751").Value = "Femacal de La Calera"
$ws.Range("C751").Value = "Coquimbo"
$ws.Range("D751").Value = 45212
$ws.Range("E751").Value = 5
$ws.Range("F751").Value = 100112037
$ws.Range("G751").Value = "Cebollín"
$ws.Range("H751").Value = "Sin especificar"
$ws.Range("I751").Value = "Primera"
$ws.Range("J751").Value = 130
$ws.Range("K751").Value = 4000
$ws.Range("L751").Value = 4000
$ws.Range("M751").Value = 4000
$ws.Range("N751").Value = "$/paquete 36 unidades"
$ws.Range("O751").Value = "Provincia de Quillota"
$ws.Range("P751").Value = 111
$ws.Range("Q751").Value = 36
$ws.Range("R751").Value = "Hortaliza"
